# Horarios actualizados Linea 141 - 214
# Refresh the scraped-schedule snapshot: update the "last updated" / "total
# rows" header cells on each sheet, patch the rows whose scrape timestamp
# moved (rows get re-labelled as later scrape passes land), and append the
# newly scraped rows at the bottom of each table.

$wb = $excel.ActiveWorkbook

# ---- Sheet "LP1912" (Linea 141, LP1912) - rows grow from 265 to 273 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2,1).Value = "Última actualización: 12:59:34"
$ws.Cells.Item(3,1).Value = "Total filas: 268"
$ws.Cells.Item(189,1).Value = "11:23:54"
$ws.Cells.Item(189,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(189,4).Value = 28
$ws.Cells.Item(190,1).Value = "10:28:12"
$ws.Cells.Item(190,3).Value = "10_OLMOS"
$ws.Cells.Item(190,4).Value = 83
$ws.Cells.Item(191,1).Value = "10:57:58"
$ws.Cells.Item(191,3).Value = "15_ABASTO"
$ws.Cells.Item(191,4).Value = 54
$ws.Cells.Item(192,1).Value = "10:28:12"
$ws.Cells.Item(192,3).Value = "215B_EL PATO"
$ws.Cells.Item(192,4).Value = 83
$ws.Cells.Item(206,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(207,3).Value = "10_OLMOS"
$ws.Cells.Item(214,1).Value = "12:16:51"
$ws.Cells.Item(214,3).Value = "16_SANTA ANA"
$ws.Cells.Item(214,4).Value = 5
$ws.Cells.Item(215,1).Value = "11:51:05"
$ws.Cells.Item(215,3).Value = "14_ABASTO"
$ws.Cells.Item(215,4).Value = 30
$ws.Cells.Item(216,1).Value = "10:28:12"
$ws.Cells.Item(216,3).Value = "215A_EL PATO"
$ws.Cells.Item(216,4).Value = 113
$ws.Cells.Item(217,1).Value = "10:28:12"
$ws.Cells.Item(217,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(217,4).Value = 113
$ws.Cells.Item(222,1).Value = "12:16:51"
$ws.Cells.Item(222,3).Value = "16_SANTA ANA"
$ws.Cells.Item(222,4).Value = 20
$ws.Cells.Item(223,1).Value = "10:57:58"
$ws.Cells.Item(223,3).Value = "27_EL RETIRO"
$ws.Cells.Item(223,4).Value = 99
$ws.Cells.Item(224,1).Value = "11:51:05"
$ws.Cells.Item(224,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(224,4).Value = 46
$ws.Cells.Item(226,1).Value = "10:57:58"
$ws.Cells.Item(226,3).Value = "17_179 Y 38"
$ws.Cells.Item(226,4).Value = 100
$ws.Cells.Item(237,1).Value = "12:59:34"
$ws.Cells.Item(237,2).Value = "12:59"
$ws.Cells.Item(237,3).Value = "16_SANTA ANA"
$ws.Cells.Item(237,4).Value = 0
$ws.Cells.Item(238,1).Value = "12:59:34"
$ws.Cells.Item(238,2).Value = "13:00"
$ws.Cells.Item(238,3).Value = "16_SANTA ANA"
$ws.Cells.Item(238,4).Value = 1
$ws.Cells.Item(239,1).Value = "12:44:21"
$ws.Cells.Item(239,2).Value = "13:02"
$ws.Cells.Item(239,3).Value = "14_ABASTO"
$ws.Cells.Item(239,4).Value = 18
$ws.Cells.Item(240,2).Value = "13:02"
$ws.Cells.Item(240,3).Value = "15_ABASTO"
$ws.Cells.Item(240,4).Value = 71
$ws.Cells.Item(241,1).Value = "12:59:34"
$ws.Cells.Item(241,2).Value = "13:04"
$ws.Cells.Item(241,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(241,4).Value = 5
$ws.Cells.Item(242,2).Value = "13:06"
$ws.Cells.Item(242,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(242,4).Value = 103
$ws.Cells.Item(243,1).Value = "11:51:05"
$ws.Cells.Item(243,2).Value = "13:07"
$ws.Cells.Item(243,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(243,4).Value = 76
$ws.Cells.Item(244,1).Value = "12:16:51"
$ws.Cells.Item(244,2).Value = "13:08"
$ws.Cells.Item(244,3).Value = "10_OLMOS"
$ws.Cells.Item(244,4).Value = 52
$ws.Cells.Item(245,2).Value = "13:13"
$ws.Cells.Item(245,3).Value = "215D_EL PATO"
$ws.Cells.Item(245,4).Value = 110
$ws.Cells.Item(246,1).Value = "12:44:21"
$ws.Cells.Item(246,2).Value = "13:14"
$ws.Cells.Item(246,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(246,4).Value = 30
$ws.Cells.Item(247,1).Value = "11:51:05"
$ws.Cells.Item(247,2).Value = "13:14"
$ws.Cells.Item(247,3).Value = "215D_EL PATO"
$ws.Cells.Item(247,4).Value = 83
$ws.Cells.Item(248,1).Value = "11:23:54"
$ws.Cells.Item(248,2).Value = "13:19"
$ws.Cells.Item(248,4).Value = 116
$ws.Cells.Item(249,2).Value = "13:20"
$ws.Cells.Item(249,3).Value = "10_OLMOS"
$ws.Cells.Item(249,4).Value = 89
$ws.Cells.Item(250,1).Value = "11:23:54"
$ws.Cells.Item(250,2).Value = "13:20"
$ws.Cells.Item(250,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(250,4).Value = 117
$ws.Cells.Item(251,1).Value = "12:44:21"
$ws.Cells.Item(251,2).Value = "13:21"
$ws.Cells.Item(251,3).Value = "10_OLMOS"
$ws.Cells.Item(251,4).Value = 37
$ws.Cells.Item(252,1).Value = "11:51:05"
$ws.Cells.Item(252,2).Value = "13:21"
$ws.Cells.Item(252,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(252,4).Value = 90
$ws.Cells.Item(253,2).Value = "13:26"
$ws.Cells.Item(253,3).Value = "14_ABASTO"
$ws.Cells.Item(253,4).Value = 70
$ws.Cells.Item(254,1).Value = "11:51:05"
$ws.Cells.Item(254,2).Value = "13:27"
$ws.Cells.Item(254,3).Value = "14_ABASTO"
$ws.Cells.Item(254,4).Value = 96
$ws.Cells.Item(255,1).Value = "12:16:51"
$ws.Cells.Item(255,2).Value = "13:32"
$ws.Cells.Item(255,3).Value = "10_OLMOS"
$ws.Cells.Item(255,4).Value = 76
$ws.Cells.Item(256,1).Value = "12:16:51"
$ws.Cells.Item(256,2).Value = "13:34"
$ws.Cells.Item(256,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(256,4).Value = 78
$ws.Cells.Item(257,1).Value = "12:44:21"
$ws.Cells.Item(257,2).Value = "13:35"
$ws.Cells.Item(257,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(257,4).Value = 51
$ws.Cells.Item(258,1).Value = "11:51:05"
$ws.Cells.Item(258,2).Value = "13:36"
$ws.Cells.Item(258,3).Value = "15_ABASTO"
$ws.Cells.Item(258,4).Value = 105
$ws.Cells.Item(259,1).Value = "11:51:05"
$ws.Cells.Item(259,2).Value = "13:46"
$ws.Cells.Item(259,3).Value = "17_ROMERO"
$ws.Cells.Item(259,4).Value = 115
$ws.Cells.Item(260,2).Value = "13:50"
$ws.Cells.Item(260,3).Value = "215A_EL PATO"
$ws.Cells.Item(260,4).Value = 94
$ws.Cells.Item(261,1).Value = "12:59:34"
$ws.Cells.Item(261,2).Value = "13:50"
$ws.Cells.Item(261,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(261,4).Value = 51
$ws.Cells.Item(262,1).Value = "12:16:51"
$ws.Cells.Item(262,2).Value = "13:55"
$ws.Cells.Item(262,3).Value = "225_GOMEZ"
$ws.Cells.Item(262,4).Value = 99
$ws.Cells.Item(263,2).Value = "13:56"
$ws.Cells.Item(263,3).Value = "225_GOMEZ"
$ws.Cells.Item(263,4).Value = 72
$ws.Cells.Item(264,1).Value = "12:16:51"
$ws.Cells.Item(264,2).Value = "14:04"
$ws.Cells.Item(264,3).Value = "17_ROMERO"
$ws.Cells.Item(264,4).Value = 108
$ws.Cells.Item(265,2).Value = "14:05"
$ws.Cells.Item(265,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(265,4).Value = 81
$ws.Cells.Item(266,1).Value = "12:44:21"
$ws.Cells.Item(266,2).Value = "14:13"
$ws.Cells.Item(266,3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(266,4).Value = 89
$ws.Cells.Item(266,5).Value = "LP1912"
$ws.Cells.Item(267,1).Value = "12:59:34"
$ws.Cells.Item(267,2).Value = "14:16"
$ws.Cells.Item(267,3).Value = "27_EL RETIRO"
$ws.Cells.Item(267,4).Value = 77
$ws.Cells.Item(267,5).Value = "LP1912"
$ws.Cells.Item(268,1).Value = "12:44:21"
$ws.Cells.Item(268,2).Value = "14:17"
$ws.Cells.Item(268,3).Value = "27_EL RETIRO"
$ws.Cells.Item(268,4).Value = 93
$ws.Cells.Item(268,5).Value = "LP1912"
$ws.Cells.Item(269,1).Value = "12:44:21"
$ws.Cells.Item(269,2).Value = "14:20"
$ws.Cells.Item(269,3).Value = "215C_EL PATO"
$ws.Cells.Item(269,4).Value = 96
$ws.Cells.Item(269,5).Value = "LP1912"
$ws.Cells.Item(270,1).Value = "12:44:21"
$ws.Cells.Item(270,2).Value = "14:21"
$ws.Cells.Item(270,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(270,4).Value = 97
$ws.Cells.Item(270,5).Value = "LP1912"
$ws.Cells.Item(271,1).Value = "12:59:34"
$ws.Cells.Item(271,2).Value = "14:44"
$ws.Cells.Item(271,3).Value = "14_ABASTO"
$ws.Cells.Item(271,4).Value = 105
$ws.Cells.Item(271,5).Value = "LP1912"
$ws.Cells.Item(272,1).Value = "12:59:34"
$ws.Cells.Item(272,2).Value = "14:56"
$ws.Cells.Item(272,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(272,4).Value = 117
$ws.Cells.Item(272,5).Value = "LP1912"
$ws.Cells.Item(273,1).Value = "12:59:34"
$ws.Cells.Item(273,2).Value = "14:58"
$ws.Cells.Item(273,3).Value = "215B_EL PATO"
$ws.Cells.Item(273,4).Value = 119
$ws.Cells.Item(273,5).Value = "LP1912"

# ---- Sheet "LP1912-215" - rows grow from 34 to 35 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2,1).Value = "Última actualización: 12:59:34"
$ws.Cells.Item(3,1).Value = "Total filas: 30"
$ws.Cells.Item(35,1).Value = "12:59:34"
$ws.Cells.Item(35,2).Value = "14:58"
$ws.Cells.Item(35,3).Value = "215B_EL PATO"
$ws.Cells.Item(35,4).Value = 119
$ws.Cells.Item(35,5).Value = "LP1912"

# ---- Sheet "6203-6173" - rows grow from 44 to 45 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2,1).Value = "Última actualización: 12:59:34"
$ws.Cells.Item(3,1).Value = "Total filas: 40"
$ws.Cells.Item(45,1).Value = "12:59:34"
$ws.Cells.Item(45,2).Value = "14:53"
$ws.Cells.Item(45,3).Value = "215D_LA PLATA"
$ws.Cells.Item(45,4).Value = 114
$ws.Cells.Item(45,5).Value = "L6203"

$wb.Save()
